$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 279.875
$ws.Range("I9").Value = 279.875
$ws.Range("K9").Value = 279.875
$ws.Range("M9").Value = -110.875

$ws.Range("H20").Value = 6455.25
$ws.Range("J20").Value = 500
$ws.Range("L20").Value = 500
$ws.Range("N20").Value = -960

$ws.Range("H35").Value = 6455.25
$ws.Range("J35").Value = 500
$ws.Range("L35").Value = 500
$ws.Range("N35").Value = -1258

$ws.Range("H38").Value = 1244.8
$ws.Range("I38").Value = 1408
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 4224
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -3852
$ws.Range("N38").Value = -3744

$ws.Range("H62").Value = 14639.5
$ws.Range("I62").Value = 2799.6667
$ws.Range("J62").Value = 19713.715
$ws.Range("K62").Value = 2799.6667
$ws.Range("L62").Value = 19713.715
$ws.Range("M62").Value = -2175.6667
$ws.Range("N62").Value = -20961.715

$ws.Range("H65").Value = 14639.5
$ws.Range("I65").Value = 2799.6667
$ws.Range("J65").Value = 19713.715
$ws.Range("K65").Value = 13998.3335
$ws.Range("L65").Value = 98568.575
$ws.Range("M65").Value = -10878.3335
$ws.Range("N65").Value = -104808.575

$ws.Range("H70").Value = 1477.0714
$ws.Range("I70").Value = 1830.8889
$ws.Range("J70").Value = 1309.4736
$ws.Range("K70").Value = 5492.6667
$ws.Range("L70").Value = 3928.4208
$ws.Range("M70").Value = -5222.6667
$ws.Range("N70").Value = -4468.4208

$ws.Range("H73").Value = 1477.0714
$ws.Range("I73").Value = 1830.8889
$ws.Range("J73").Value = 1309.4736
$ws.Range("K73").Value = 5492.6667
$ws.Range("L73").Value = 3928.4208
$ws.Range("M73").Value = -4556.6667
$ws.Range("N73").Value = -5800.4208

$ws.Range("H101").Value = 3051.7144
$ws.Range("I101").Value = 3029.5
$ws.Range("J101").Value = 3185
$ws.Range("K101").Value = 9088.5
$ws.Range("L101").Value = 9555
$ws.Range("M101").Value = -7466.5
$ws.Range("N101").Value = -12799

$ws.Range("H107").Value = 383.42856
$ws.Range("I107").Value = 355.66666
$ws.Range("K107").Value = 355.66666
$ws.Range("M107").Value = 1564.33334

$ws.Range("H132").Value = 2528.6667
$ws.Range("I132").Value = 2602.7827
$ws.Range("K132").Value = 7808.348100000001
$ws.Range("M132").Value = -5278.348100000001

$ws.Range("H138").Value = 3163
$ws.Range("I138").Value = 1182.4286
$ws.Range("K138").Value = 3547.2858
$ws.Range("M138").Value = 1592.7142

$ws.Range("H141").Value = 5019
$ws.Range("I141").Value = 5019
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 15057
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -9877
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10419169
$ws.Range("I32").Value = 11906359
$ws.Range("K32").Value = 11906359
$ws.Range("M32").Value = -11906072

$ws.Range("H38").Value = 39750
$ws.Range("J38").Value = 39750
$ws.Range("L38").Value = 39750
$ws.Range("N38").Value = -40684

$ws.Range("H45").Value = 2270.6667
$ws.Range("I45").Value = 1869.091
$ws.Range("K45").Value = 1869.091
$ws.Range("M45").Value = -1492.091

$ws.Range("H124").Value = 80806.336
$ws.Range("J124").Value = 80806.336
$ws.Range("L124").Value = 80806.336
$ws.Range("N124").Value = -90626.336

$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

$ws.Range("H132").Value = 8221.75
$ws.Range("I132").Value = 2897.0908
$ws.Range("K132").Value = 8691.2724
$ws.Range("M132").Value = -6161.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2310
$ws.Range("I105").Value = 1275
$ws.Range("K105").Value = 1275
$ws.Range("M105").Value = 472

$ws.Range("H107").Value = 1927.7
$ws.Range("I107").Value = 1310.9286
$ws.Range("K107").Value = 1310.9286
$ws.Range("M107").Value = 609.0714

$ws.Range("H132").Value = 95294.12
$ws.Range("I132").Value = 96857.14
$ws.Range("K132").Value = 96857.14
$ws.Range("M132").Value = -91797.14

$ws.Range("H134").Value = 82985.62
$ws.Range("I134").Value = 1878.5
$ws.Range("K134").Value = 5635.5
$ws.Range("M134").Value = -3100.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1052
$ws.Range("I16").Value = 1052
$ws.Range("K16").Value = 1052
$ws.Range("M16").Value = -765

$ws.Range("H31").Value = 438217.66
$ws.Range("I31").Value = 5397.0884
$ws.Range("K31").Value = 5397.0884
$ws.Range("M31").Value = -5102.0884

$ws.Range("H34").Value = 438217.66
$ws.Range("I34").Value = 5397.0884
$ws.Range("K34").Value = 5397.0884
$ws.Range("M34").Value = -5195.0884

$ws.Range("H86").Value = 3982.1667
$ws.Range("I86").Value = 3982.1667
$ws.Range("K86").Value = 3982.1667
$ws.Range("M86").Value = -2859.1667

$ws.Range("H89").Value = 3982.1667
$ws.Range("I89").Value = 3982.1667
$ws.Range("K89").Value = 19910.8335
$ws.Range("M89").Value = -14294.8335

$ws.Range("H105").Value = 1930.1666
$ws.Range("I105").Value = 2071.182
$ws.Range("J105").Value = 379
$ws.Range("K105").Value = 2071.182
$ws.Range("L105").Value = 379
$ws.Range("M105").Value = -324.1819999999998
$ws.Range("N105").Value = -3873

$ws.Range("H113").Value = 1052
$ws.Range("I113").Value = 1052
$ws.Range("K113").Value = 1052
$ws.Range("M113").Value = 1118

$ws.Range("H132").Value = 2156.9
$ws.Range("I132").Value = 1954.6316
$ws.Range("K132").Value = 5863.8948
$ws.Range("M132").Value = -3333.8948

$ws.Range("H134").Value = 456096.47
$ws.Range("I134").Value = 456096.47
$ws.Range("K134").Value = 1368289.41
$ws.Range("M134").Value = -1365754.41

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 4111
$ws.Range("J121").Value = 4412.5713
$ws.Range("L121").Value = 13237.7139
$ws.Range("N121").Value = -15857.7139

$ws.Range("H131").Value = 11449.5
$ws.Range("I131").Value = 18633
$ws.Range("J131").Value = 7139.4
$ws.Range("K131").Value = 55899
$ws.Range("L131").Value = 21418.2
$ws.Range("M131").Value = -50859
$ws.Range("N131").Value = -31498.2

$ws.Range("H141").Value = 281231.28
$ws.Range("J141").Value = 11071.286
$ws.Range("L141").Value = 33213.858
$ws.Range("N141").Value = -43573.858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 4999
$ws.Range("I41").Value = 4998
$ws.Range("K41").Value = 4998
$ws.Range("M41").Value = -4643

$ws.Range("H80").Value = 3379.389
$ws.Range("I80").Value = 3239.3125
$ws.Range("K80").Value = 3239.3125
$ws.Range("M80").Value = -2241.3125

$ws.Range("H83").Value = 3379.389
$ws.Range("I83").Value = 3239.3125
$ws.Range("K83").Value = 16196.5625
$ws.Range("M83").Value = -11204.5625

$ws.Range("H122").Value = 1162.1111
$ws.Range("I122").Value = 933.9231
$ws.Range("K122").Value = 2801.7693
$ws.Range("M122").Value = -351.7692999999999

$ws.Range("H132").Value = 50002590
$ws.Range("I132").Value = 71431580
$ws.Range("J132").Value = 1618.6666
$ws.Range("K132").Value = 214294740
$ws.Range("L132").Value = 4855.9998
$ws.Range("M132").Value = -214292210
$ws.Range("N132").Value = -9915.9998

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 45454990
$ws.Range("I55").Value = 58823960
$ws.Range("J55").Value = 502.6
$ws.Range("K55").Value = 58823960
$ws.Range("L55").Value = 502.6
$ws.Range("M55").Value = -58823787
$ws.Range("N55").Value = -848.6

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H132").Value = 40246.594
$ws.Range("I132").Value = 3679.6316
$ws.Range("K132").Value = 11038.8948
$ws.Range("M132").Value = -8508.8948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 114000
$ws.Range("I2").Value = 130000
$ws.Range("K2").Value = 130000
$ws.Range("M2").Value = -129888

$ws.Range("H100").Value = 1420.5555
$ws.Range("I100").Value = 1448.125
$ws.Range("J100").Value = 1200
$ws.Range("K100").Value = 2896.25
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -2355.25
$ws.Range("N100").Value = -3482

$ws.Range("H113").Value = 2149.5
$ws.Range("I113").Value = 4666.3335
$ws.Range("J113").Value = 639.4
$ws.Range("K113").Value = 13999.0005
$ws.Range("L113").Value = 1918.2
$ws.Range("M113").Value = -11829.0005
$ws.Range("N113").Value = -6258.2

$ws.Range("H132").Value = 1604.2222
$ws.Range("I132").Value = 1576.2941
$ws.Range("K132").Value = 4728.8823
$ws.Range("L132").Value = 5063.6362
$ws.Range("M132").Value = -2198.8823
